$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-annotated dialog act (DAMSLTag / DialogAct) values for specific rows
# following clean-up work to the original transcripts (re-run of SGNN tagger).
$updates = @(
    @{ Row = 27; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 28; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 30; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 32; DamslTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 38; DamslTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 41; DamslTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 42; DamslTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 47; DamslTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 50; DamslTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 52; DamslTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 56; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 58; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 76; DamslTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 77; DamslTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 83; DamslTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 84; DamslTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 90; DamslTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 91; DamslTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 95; DamslTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 104; DamslTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 120; DamslTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 132; DamslTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 140; DamslTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 154; DamslTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 156; DamslTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 170; DamslTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 191; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 192; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 216; DamslTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 220; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 225; DamslTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 230; DamslTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 236; DamslTag = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 237; DamslTag = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 239; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 245; DamslTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 249; DamslTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 252; DamslTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 255; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 260; DamslTag = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 269; DamslTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 272; DamslTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 273; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 274; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 281; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 284; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 285; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 305; DamslTag = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 312; DamslTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 322; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 324; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 328; DamslTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 331; DamslTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 335; DamslTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 343; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 346; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 353; DamslTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 354; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 372; DamslTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 389; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 393; DamslTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 414; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 435; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 442; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 447; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 448; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 450; DamslTag = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 454; DamslTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 460; DamslTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 469; DamslTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 470; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 471; DamslTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 480; DamslTag = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 487; DamslTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 490; DamslTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 500; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 503; DamslTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 504; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 506; DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 534; DamslTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 535; DamslTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 543; DamslTag = "qy"; DialogAct = "Yes-No-Question" }
    @{ Row = 552; DamslTag = "ba"; DialogAct = "Appreciation" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.DamslTag
    $ws.Cells.Item($u.Row, 10).Value = $u.DialogAct
}

Write-Output "Updated $($updates.Count) rows"